$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3048.75
$ws.Range("J40").Value = 3048.75
$ws.Range("L40").Value = 3048.75
$ws.Range("N40").Value = -3398.75

$ws.Range("H43").Value = 800
$ws.Range("I43").Value = 800
$ws.Range("K43").Value = 800
$ws.Range("M43").Value = -731

$ws.Range("H70").Value = 2993.7856
$ws.Range("I70").Value = 1319
$ws.Range("J70").Value = 4249.875
$ws.Range("K70").Value = 3957
$ws.Range("L70").Value = 12749.625
$ws.Range("M70").Value = -3687
$ws.Range("N70").Value = -13289.625

$ws.Range("H73").Value = 2993.7856
$ws.Range("I73").Value = 1319
$ws.Range("J73").Value = 4249.875
$ws.Range("K73").Value = 3957
$ws.Range("L73").Value = 12749.625
$ws.Range("M73").Value = -3021
$ws.Range("N73").Value = -14621.625

$ws.Range("H99").Value = 1301.8572
$ws.Range("I99").Value = 1080.75
$ws.Range("J99").Value = 1596.6666
$ws.Range("K99").Value = 3242.25
$ws.Range("L99").Value = 4789.9998
$ws.Range("M99").Value = -1744.25
$ws.Range("N99").Value = -7785.9998

$ws.Range("H112").Value = 1993
$ws.Range("J112").Value = 2332.8572
$ws.Range("L112").Value = 6998.571599999999
$ws.Range("N112").Value = -9214.571599999999

$ws.Range("H135").Value = 560.3333
$ws.Range("I135").Value = 190.5
$ws.Range("J135").Value = 1300
$ws.Range("K135").Value = 1714.5
$ws.Range("L135").Value = 11700
$ws.Range("M135").Value = 820.5
$ws.Range("N135").Value = -16770

$ws.Range("H138").Value = 2223.2
$ws.Range("I138").Value = 1529
$ws.Range("J138").Value = 5000
$ws.Range("K138").Value = 4587
$ws.Range("L138").Value = 15000
$ws.Range("M138").Value = 553
$ws.Range("N138").Value = -25280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1268.8462
$ws.Range("I2").Value = 724.6667
$ws.Range("K2").Value = 724.6667
$ws.Range("M2").Value = -611.6667

$ws.Range("H32").Value = 4507.4
$ws.Range("I32").Value = 4404.207
$ws.Range("J32").Value = 7500
$ws.Range("K32").Value = 4404.207
$ws.Range("L32").Value = 7500
$ws.Range("M32").Value = -4117.207
$ws.Range("N32").Value = -8074

$ws.Range("H45").Value = 1918.3334
$ws.Range("I45").Value = 1834.375
$ws.Range("J45").Value = 2014.2858
$ws.Range("K45").Value = 1834.375
$ws.Range("L45").Value = 2014.2858
$ws.Range("M45").Value = -1457.375
$ws.Range("N45").Value = -2768.2858

$ws.Range("H61").Value = 3498.5
$ws.Range("I61").Value = 3498
$ws.Range("K61").Value = 3498
$ws.Range("M61").Value = -3286

$ws.Range("H63").Value = 2850
$ws.Range("I63").Value = 2390
$ws.Range("K63").Value = 2390
$ws.Range("M63").Value = -1704

$ws.Range("H66").Value = 2850
$ws.Range("I66").Value = 2390
$ws.Range("K66").Value = 11950
$ws.Range("M66").Value = -8518

$ws.Range("H97").Value = 696
$ws.Range("I97").Value = 655.73334
$ws.Range("K97").Value = 655.73334
$ws.Range("M97").Value = -159.73334

$ws.Range("H101").Value = 80066.664
$ws.Range("J101").Value = 80066.664
$ws.Range("L101").Value = 80066.664
$ws.Range("N101").Value = -86556.664

$ws.Range("H102").Value = 5171.25
$ws.Range("I102").Value = 1845
$ws.Range("J102").Value = 8497.5
$ws.Range("K102").Value = 1845
$ws.Range("L102").Value = 8497.5
$ws.Range("M102").Value = -223
$ws.Range("N102").Value = -11741.5

$ws.Range("H105").Value = 55000
$ws.Range("J105").Value = 55000
$ws.Range("L105").Value = 55000
$ws.Range("N105").Value = -61988

$ws.Range("H110").Value = 2575.875
$ws.Range("J110").Value = 2500
$ws.Range("L110").Value = 2500
$ws.Range("N110").Value = -6590

$ws.Range("H116").Value = 1268.8462
$ws.Range("I116").Value = 724.6667
$ws.Range("K116").Value = 724.6667
$ws.Range("M116").Value = 1569.3333

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H118").Value = 69999
$ws.Range("J118").Value = 69999
$ws.Range("L118").Value = 69999
$ws.Range("N118").Value = -73313

$ws.Range("H123").Value = 80000
$ws.Range("J123").Value = 80000
$ws.Range("L123").Value = 80000
$ws.Range("N123").Value = -89800

$ws.Range("H132").Value = 2999
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2999
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 8997
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -14057

$ws.Range("H136").Value = 3498.5
$ws.Range("I136").Value = 3498
$ws.Range("K136").Value = 10494
$ws.Range("M136").Value = -7944

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1268.8462
$ws.Range("I3").Value = 724.6667
$ws.Range("K3").Value = 724.6667
$ws.Range("M3").Value = -610.6667

$ws.Range("H105").Value = 1883.3334
$ws.Range("I105").Value = 1575
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 1575
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = 172
$ws.Range("N105").Value = -5994

$ws.Range("H134").Value = 4711.5557
$ws.Range("I134").Value = 4758.154
$ws.Range("K134").Value = 14274.462
$ws.Range("M134").Value = -11739.462

$ws.Range("H135").Value = 49999
$ws.Range("J135").Value = 49999
$ws.Range("L135").Value = 49999
$ws.Range("N135").Value = -60139

$ws.Range("H137").Value = 34999.5
$ws.Range("I137").Value = 34999.5
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 34999.5
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -29899.5
$ws.Range("N137").ClearContents()

$ws.Range("H140").Value = 80000
$ws.Range("J140").Value = 80000
$ws.Range("L140").Value = 80000
$ws.Range("N140").Value = -90360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 18400

$ws.Range("H93").Value = 17499.5
$ws.Range("I93").Value = 17499.5
$ws.Range("K93").Value = 17499.5
$ws.Range("M93").Value = -15627.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 999.25
$ws.Range("J5").Value = 999
$ws.Range("L5").Value = 2997
$ws.Range("N5").Value = -3221

$ws.Range("H34").Value = 558.2
$ws.Range("I34").Value = 456.8
$ws.Range("J34").Value = 659.6
$ws.Range("K34").Value = 1370.4
$ws.Range("L34").Value = 1978.8
$ws.Range("M34").Value = -1286.4
$ws.Range("N34").Value = -2146.8

$ws.Range("H39").Value = 2466.6667
$ws.Range("J39").Value = 2700
$ws.Range("L39").Value = 8100
$ws.Range("N39").Value = -8688

$ws.Range("H55").Value = 2714.2856
$ws.Range("J55").Value = 5000
$ws.Range("L55").Value = 15000
$ws.Range("N55").Value = -15354

$ws.Range("H135").Value = 999.25
$ws.Range("J135").Value = 999
$ws.Range("L135").Value = 8991
$ws.Range("N135").Value = -14061

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5184.75
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 5184.75
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 5184.75
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -7180.75

$ws.Range("H83").Value = 5184.75
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 5184.75
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 25923.75
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -35907.75

$ws.Range("H113").Value = 4969.125
$ws.Range("J113").Value = 5998
$ws.Range("L113").Value = 5998
$ws.Range("N113").Value = -10338

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1500
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 1500
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 1500
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -2222

$ws.Range("H85").Value = 1500
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 1500
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 1500
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -3996

$ws.Range("H100").Value = 1263.3334
$ws.Range("I100").Value = 1263.3334
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1263.3334
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -722.3334
$ws.Range("N100").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 10000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 10000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 20000
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -22122

$ws.Range("H84").Value = 10000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 10000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 100000
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -110608

$ws.Range("H96").Value = 804
$ws.Range("I96").Value = 804
$ws.Range("K96").Value = 804
$ws.Range("M96").Value = 569

$ws.Range("H132").Value = 2030.8125
$ws.Range("I132").Value = 1883.6666
$ws.Range("K132").Value = 5650.9998
$ws.Range("M132").Value = -3120.9998

$ws.Range("H136").Value = 2009.6786
$ws.Range("I136").Value = 1406.88
$ws.Range("K136").Value = 4220.64
$ws.Range("M136").Value = -1670.64
